$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Check Sinh viên"
$ws.Range("B16").Value = "Kiểm tra link Về TDB"
$ws.Range("C16").Value = "WebElement: https://tbd.edu.vn/sinh-vien/"
$ws.Range("D16").Value = "Test Link Passed!"
$ws.Range("E16").Value = "Sinh viên"
$ws.Range("F16").Value = "Pass"
$ws.Range("G16").Value = "Phạm Minh Tuấn"

$ws.Range("G16").Select()
